$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "60 Hz [Hz]"
$ws.Range("H5").Value = "10 kHz [kHz]"
$ws.Range("I5").Value = "100 kHz [kHz]"

$ws.Range("F6").Value = "0,1 s"
$ws.Range("G6").Value = 70
$ws.Range("H6").Value = 10.63
$ws.Range("I6").Value = 110.83

$ws.Range("F7").Value = "1 s"
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = 10.622
$ws.Range("I7").Value = 110.828

$ws.Range("F8").Value = "10 s"
$ws.Range("G8").Value = 67.9
$ws.Range("H8").Value = 10.6231
$ws.Range("I8").Value = 110.8373

$ws.Range("G10").Value = "60 Hz [s]"
$ws.Range("H10").Value = "10 kHz [s]"
$ws.Range("I10").Value = "100 kHz [s]"

$ws.Range("F11").Value = "1 T"
$ws.Range("G11").Value = 15100.7
$ws.Range("H11").Value = 93.6
$ws.Range("I11").Value = 9

$ws.Range("F12").Value = "10 T "
$ws.Range("G12").Value = 15101.87
$ws.Range("H12").Value = 93.64
$ws.Range("I12").Value = "-"

$ws.Range("F13").Value = "100 T"
$ws.Range("G13").Value = "-"
$ws.Range("H13").Value = 93.6488
$ws.Range("I13").Value = "-"

$ws.Range("F14").Value = "1000 T"
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = 93.63633
$ws.Range("I14").Value = 9.05553

$ws.Range("F15").Value = " "

$ws.Range("G18").Value = "f [kHz]"
$ws.Range("H18").Value = "T [us]"

$ws.Range("G19").Value = 32.7677
$ws.Range("H19").Value = 2000038.9

$ws.Range("K5").Select()
